$d = $word.ActiveDocument

# Replace each equation cell text with its new value.
# wdFindContinue = 1, wdReplaceAll = 2
$d.Content.Find.Execute("59+14=73", $true, $false, $false, $false, $false, $true, 1, $false, "11+37=48", 2) | Out-Null
$d.Content.Find.Execute("20+15=35", $true, $false, $false, $false, $false, $true, 1, $false, "91-31=60", 2) | Out-Null
$d.Content.Find.Execute("95-2=93", $true, $false, $false, $false, $false, $true, 1, $false, "56+15=71", 2) | Out-Null
$d.Content.Find.Execute("56-21=35", $true, $false, $false, $false, $false, $true, 1, $false, "28+31=59", 2) | Out-Null
$d.Content.Find.Execute("99-26=73", $true, $false, $false, $false, $false, $true, 1, $false, "44-1=43", 2) | Out-Null
$d.Content.Find.Execute("85-30=55", $true, $false, $false, $false, $false, $true, 1, $false, "94-39=55", 2) | Out-Null
$d.Content.Find.Execute("98-35=63", $true, $false, $false, $false, $false, $true, 1, $false, "55-53=2", 2) | Out-Null
$d.Content.Find.Execute("8+22=30", $true, $false, $false, $false, $false, $true, 1, $false, "71+26=97", 2) | Out-Null
$d.Content.Find.Execute("61-41=20", $true, $false, $false, $false, $false, $true, 1, $false, "2+41=43", 2) | Out-Null
$d.Content.Find.Execute("52+19=71", $true, $false, $false, $false, $false, $true, 1, $false, "69-57=12", 2) | Out-Null
$d.Content.Find.Execute("34+29=63", $true, $false, $false, $false, $false, $true, 1, $false, "55+18=73", 2) | Out-Null
$d.Content.Find.Execute("3+20=23", $true, $false, $false, $false, $false, $true, 1, $false, "30+23=53", 2) | Out-Null
$d.Content.Find.Execute("82+8=90", $true, $false, $false, $false, $false, $true, 1, $false, "44+48=92", 2) | Out-Null
$d.Content.Find.Execute("38-24=14", $true, $false, $false, $false, $false, $true, 1, $false, "44+26=70", 2) | Out-Null
$d.Content.Find.Execute("56-16=40", $true, $false, $false, $false, $false, $true, 1, $false, "27-0=27", 2) | Out-Null
$d.Content.Find.Execute("94-12=82", $true, $false, $false, $false, $false, $true, 1, $false, "49+11=60", 2) | Out-Null
$d.Content.Find.Execute("69-22=47", $true, $false, $false, $false, $false, $true, 1, $false, "19-14=5", 2) | Out-Null
$d.Content.Find.Execute("35+38=73", $true, $false, $false, $false, $false, $true, 1, $false, "93-78=15", 2) | Out-Null
$d.Content.Find.Execute("76-25=51", $true, $false, $false, $false, $false, $true, 1, $false, "60+38=98", 2) | Out-Null
$d.Content.Find.Execute("8+5=13", $true, $false, $false, $false, $false, $true, 1, $false, "88-80=8", 2) | Out-Null
$d.Content.Find.Execute("26-23=3", $true, $false, $false, $false, $false, $true, 1, $false, "24+18=42", 2) | Out-Null
$d.Content.Find.Execute("59-55=4", $true, $false, $false, $false, $false, $true, 1, $false, "92-49=43", 2) | Out-Null
$d.Content.Find.Execute("48-32=16", $true, $false, $false, $false, $false, $true, 1, $false, "74-55=19", 2) | Out-Null
$d.Content.Find.Execute("65-52=13", $true, $false, $false, $false, $false, $true, 1, $false, "68-9=59", 2) | Out-Null
$d.Content.Find.Execute("52-24=28", $true, $false, $false, $false, $false, $true, 1, $false, "20+14=34", 2) | Out-Null
$d.Content.Find.Execute("62-18=44", $true, $false, $false, $false, $false, $true, 1, $false, "90-83=7", 2) | Out-Null
$d.Content.Find.Execute("65-59=6", $true, $false, $false, $false, $false, $true, 1, $false, "99-87=12", 2) | Out-Null
$d.Content.Find.Execute("92-39=53", $true, $false, $false, $false, $false, $true, 1, $false, "57+36=93", 2) | Out-Null
$d.Content.Find.Execute("84-6=78", $true, $false, $false, $false, $false, $true, 1, $false, "56+38=94", 2) | Out-Null
$d.Content.Find.Execute("77-9=68", $true, $false, $false, $false, $false, $true, 1, $false, "90-56=34", 2) | Out-Null
$d.Content.Find.Execute("4+23=27", $true, $false, $false, $false, $false, $true, 1, $false, "83-77=6", 2) | Out-Null
$d.Content.Find.Execute("9+49=58", $true, $false, $false, $false, $false, $true, 1, $false, "41+45=86", 2) | Out-Null
$d.Content.Find.Execute("81-75=6", $true, $false, $false, $false, $false, $true, 1, $false, "51+38=89", 2) | Out-Null
$d.Content.Find.Execute("94-55=39", $true, $false, $false, $false, $false, $true, 1, $false, "50-29=21", 2) | Out-Null
$d.Content.Find.Execute("12+82=94", $true, $false, $false, $false, $false, $true, 1, $false, "78-41=37", 2) | Out-Null
$d.Content.Find.Execute("83+0=83", $true, $false, $false, $false, $false, $true, 1, $false, "65+18=83", 2) | Out-Null
$d.Content.Find.Execute("15-7=8", $true, $false, $false, $false, $false, $true, 1, $false, "30+30=60", 2) | Out-Null
$d.Content.Find.Execute("94-49=45", $true, $false, $false, $false, $false, $true, 1, $false, "8-3=5", 2) | Out-Null
$d.Content.Find.Execute("8+15=23", $true, $false, $false, $false, $false, $true, 1, $false, "56+2=58", 2) | Out-Null
$d.Content.Find.Execute("31+13=44", $true, $false, $false, $false, $false, $true, 1, $false, "54-51=3", 2) | Out-Null
$d.Content.Find.Execute("73-26=47", $true, $false, $false, $false, $false, $true, 1, $false, "90-70=20", 2) | Out-Null
$d.Content.Find.Execute("14+38=52", $true, $false, $false, $false, $false, $true, 1, $false, "19+60=79", 2) | Out-Null
$d.Content.Find.Execute("27+42=69", $true, $false, $false, $false, $false, $true, 1, $false, "89-88=1", 2) | Out-Null
$d.Content.Find.Execute("5+1=6", $true, $false, $false, $false, $false, $true, 1, $false, "66-25=41", 2) | Out-Null
$d.Content.Find.Execute("62-35=27", $true, $false, $false, $false, $false, $true, 1, $false, "26-1=25", 2) | Out-Null
$d.Content.Find.Execute("9+70=79", $true, $false, $false, $false, $false, $true, 1, $false, "59+20=79", 2) | Out-Null
$d.Content.Find.Execute("82-76=6", $true, $false, $false, $false, $false, $true, 1, $false, "0+66=66", 2) | Out-Null
$d.Content.Find.Execute("42+21=63", $true, $false, $false, $false, $false, $true, 1, $false, "34+9=43", 2) | Out-Null
$d.Content.Find.Execute("0+57=57", $true, $false, $false, $false, $false, $true, 1, $false, "75-15=60", 2) | Out-Null
$d.Content.Find.Execute("98-0=98", $true, $false, $false, $false, $false, $true, 1, $false, "22-7=15", 2) | Out-Null
$d.Content.Find.Execute("69-11=58", $true, $false, $false, $false, $false, $true, 1, $false, "56+6=62", 2) | Out-Null
$d.Content.Find.Execute("55+44=99", $true, $false, $false, $false, $false, $true, 1, $false, "75-8=67", 2) | Out-Null
$d.Content.Find.Execute("11-5=6", $true, $false, $false, $false, $false, $true, 1, $false, "89+4=93", 2) | Out-Null
$d.Content.Find.Execute("36-27=9", $true, $false, $false, $false, $false, $true, 1, $false, "1+7=8", 2) | Out-Null
$d.Content.Find.Execute("58-27=31", $true, $false, $false, $false, $false, $true, 1, $false, "27+54=81", 2) | Out-Null
$d.Content.Find.Execute("21+40=61", $true, $false, $false, $false, $false, $true, 1, $false, "52-42=10", 2) | Out-Null
$d.Content.Find.Execute("43+7=50", $true, $false, $false, $false, $false, $true, 1, $false, "41+19=60", 2) | Out-Null
$d.Content.Find.Execute("83-45=38", $true, $false, $false, $false, $false, $true, 1, $false, "60-50=10", 2) | Out-Null
$d.Content.Find.Execute("8+73=81", $true, $false, $false, $false, $false, $true, 1, $false, "21+14=35", 2) | Out-Null
$d.Content.Find.Execute("94-34=60", $true, $false, $false, $false, $false, $true, 1, $false, "70+21=91", 2) | Out-Null
$d.Content.Find.Execute("73-23=50", $true, $false, $false, $false, $false, $true, 1, $false, "38+23=61", 2) | Out-Null
$d.Content.Find.Execute("78-32=46", $true, $false, $false, $false, $false, $true, 1, $false, "37+9=46", 2) | Out-Null
$d.Content.Find.Execute("38+18=56", $true, $false, $false, $false, $false, $true, 1, $false, "79-9=70", 2) | Out-Null
$d.Content.Find.Execute("93-15=78", $true, $false, $false, $false, $false, $true, 1, $false, "73+0=73", 2) | Out-Null
$d.Content.Find.Execute("1+80=81", $true, $false, $false, $false, $false, $true, 1, $false, "3+15=18", 2) | Out-Null
$d.Content.Find.Execute("55+8=63", $true, $false, $false, $false, $false, $true, 1, $false, "74+22=96", 2) | Out-Null
$d.Content.Find.Execute("1+61=62", $true, $false, $false, $false, $false, $true, 1, $false, "97-21=76", 2) | Out-Null
$d.Content.Find.Execute("34+48=82", $true, $false, $false, $false, $false, $true, 1, $false, "41+26=67", 2) | Out-Null
$d.Content.Find.Execute("99-54=45", $true, $false, $false, $false, $false, $true, 1, $false, "19+34=53", 2) | Out-Null
$d.Content.Find.Execute("14+4=18", $true, $false, $false, $false, $false, $true, 1, $false, "99-51=48", 2) | Out-Null
$d.Content.Find.Execute("59-37=22", $true, $false, $false, $false, $false, $true, 1, $false, "30-21=9", 2) | Out-Null
$d.Content.Find.Execute("81-66=15", $true, $false, $false, $false, $false, $true, 1, $false, "87-86=1", 2) | Out-Null
$d.Content.Find.Execute("78+6=84", $true, $false, $false, $false, $false, $true, 1, $false, "66+2=68", 2) | Out-Null
$d.Content.Find.Execute("88-58=30", $true, $false, $false, $false, $false, $true, 1, $false, "37+17=54", 2) | Out-Null
$d.Content.Find.Execute("91-87=4", $true, $false, $false, $false, $false, $true, 1, $false, "87-67=20", 2) | Out-Null
$d.Content.Find.Execute("5+43=48", $true, $false, $false, $false, $false, $true, 1, $false, "14+49=63", 2) | Out-Null
$d.Content.Find.Execute("61-3=58", $true, $false, $false, $false, $false, $true, 1, $false, "52+9=61", 2) | Out-Null
$d.Content.Find.Execute("62-56=6", $true, $false, $false, $false, $false, $true, 1, $false, "91-76=15", 2) | Out-Null
$d.Content.Find.Execute("41+36=77", $true, $false, $false, $false, $false, $true, 1, $false, "40+20=60", 2) | Out-Null
$d.Content.Find.Execute("12+16=28", $true, $false, $false, $false, $false, $true, 1, $false, "95-17=78", 2) | Out-Null
$d.Content.Find.Execute("9+47=56", $true, $false, $false, $false, $false, $true, 1, $false, "54+2=56", 2) | Out-Null
$d.Content.Find.Execute("48-19=29", $true, $false, $false, $false, $false, $true, 1, $false, "21+71=92", 2) | Out-Null
$d.Content.Find.Execute("57+9=66", $true, $false, $false, $false, $false, $true, 1, $false, "33+13=46", 2) | Out-Null
$d.Content.Find.Execute("49-22=27", $true, $false, $false, $false, $false, $true, 1, $false, "7+30=37", 2) | Out-Null
$d.Content.Find.Execute("96-17=79", $true, $false, $false, $false, $false, $true, 1, $false, "45-18=27", 2) | Out-Null
$d.Content.Find.Execute("20-9=11", $true, $false, $false, $false, $false, $true, 1, $false, "75-8=67", 2) | Out-Null
$d.Content.Find.Execute("59-11=48", $true, $false, $false, $false, $false, $true, 1, $false, "85-27=58", 2) | Out-Null
$d.Content.Find.Execute("69-6=63", $true, $false, $false, $false, $false, $true, 1, $false, "13+83=96", 2) | Out-Null
$d.Content.Find.Execute("65+22=87", $true, $false, $false, $false, $false, $true, 1, $false, "63-43=20", 2) | Out-Null
$d.Content.Find.Execute("94-67=27", $true, $false, $false, $false, $false, $true, 1, $false, "6-3=3", 2) | Out-Null
$d.Content.Find.Execute("20+28=48", $true, $false, $false, $false, $false, $true, 1, $false, "14+59=73", 2) | Out-Null
$d.Content.Find.Execute("83-42=41", $true, $false, $false, $false, $false, $true, 1, $false, "29+48=77", 2) | Out-Null
$d.Content.Find.Execute("8+29=37", $true, $false, $false, $false, $false, $true, 1, $false, "46+22=68", 2) | Out-Null
$d.Content.Find.Execute("16+39=55", $true, $false, $false, $false, $false, $true, 1, $false, "18+67=85", 2) | Out-Null
$d.Content.Find.Execute("91-65=26", $true, $false, $false, $false, $false, $true, 1, $false, "99-86=13", 2) | Out-Null
$d.Content.Find.Execute("12+47=59", $true, $false, $false, $false, $false, $true, 1, $false, "17+4=21", 2) | Out-Null
$d.Content.Find.Execute("80-39=41", $true, $false, $false, $false, $false, $true, 1, $false, "50+45=95", 2) | Out-Null
$d.Content.Find.Execute("49+47=96", $true, $false, $false, $false, $false, $true, 1, $false, "6+45=51", 2) | Out-Null
$d.Content.Find.Execute("5+67=72", $true, $false, $false, $false, $false, $true, 1, $false, "54+12=66", 2) | Out-Null
$d.Content.Find.Execute("23+67=90", $true, $false, $false, $false, $false, $true, 1, $false, "30-26=4", 2) | Out-Null
